# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 3120e744-... handoff row on both the zh-cn and de-de
# status sheets, reflecting the newly generated report run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for 3120e744-b469-4407-8802-41b07d275d20 (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 08:39:11"   # Correspond Handoff Datetime
$wsZhCn.Range("H4").Value = "2016-03-21 08:39:32"   # Correspond Handback DateTime

# de-de sheet: row for 3120e744-b469-4407-8802-41b07d275d20 (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 08:39:14"   # Correspond Handoff Datetime
$wsDeDe.Range("H4").Value = "2016-03-21 08:39:38"   # Correspond Handback DateTime
